$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.195.48'
$ws.Range("E2").Value = '  +2.04%  '
$ws.Range("D3").Value = '1.992.10'
$ws.Range("E3").Value = '  +5.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7957'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +67.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '255.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9991'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3498'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +20.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.15'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +25.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06998'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8466'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08181'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '100.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").Value = '1.990.85'
$ws.Range("E14").Value = '  +5.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.634'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +16.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '273.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.22%  '
$ws.Range("D18").Value = '31.191.68'
$ws.Range("E18").Value = '  +2.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.888'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007955'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.49%  '
$ws.Range("D21").Value = '2.249.93'
$ws.Range("E21").Value = '  +5.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9990'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.063'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1508'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +55.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.345'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +22.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.599'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.357'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.590'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.432'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05262'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7819'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.217'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.51%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9983'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02008'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.59%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.655'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '79.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4660'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.126'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8535'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9992'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.683'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.945'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4303'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.37%  '
